$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold decimal-looking numeric text which must remain exact text
# (avoids Excel auto-converting them to floating point numbers and losing
# trailing zeros / exact formatting).
$textCells = @("D5", "D6", "D7", "D10", "D12", "D14", "D17", "D19", "D20", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D40", "D41", "D42", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '42.950.78'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.326.47'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '302.06'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').Value = '95.89'
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('D7').Value = '0.503'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').Value = '34.35'
$ws.Range('E10').Value = '  -3.08%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '0.0785'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').Value = '6.75'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '2.692.15'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = '2.354.39'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '0.791'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '42.891.38'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '12.25'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '6.17'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '67.96'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '236.20'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '2.25'
$ws.Range('E24').Value = '  +4.09%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').Value = '24.69'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('E28').Value = '  -5.95%  '
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').Value = '32.24'
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').Value = '144.61'
$ws.Range('E31').Value = '  -12.96%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').Value = '5.01'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('D34').Value = '17.88'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('D35').Value = '0.0702'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('D40').Value = '2.75'
$ws.Range('D41').Value = '22.21'
$ws.Range('E41').Value = '  +23.55%  '
$ws.Range('D42').Value = '0.108'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').Value = '1.933.10'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').Value = '10.12'
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('D46').Value = '2.06'
$ws.Range('E46').Value = '  -2.45%  '
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').Value = '2.88'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.559.52'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').Value = '53.61'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '73.01'
$ws.Range('E51').Value = '  +1.67%  '
